$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot original values for columns D, L, M, N, O, P, S across rows 4-33
$cols = @("D","L","M","N","O","P","S")
$orig = @{}
for ($r = 4; $r -le 33; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $orig[$r] = $rowVals
}

# Mapping: target row -> source row (values to copy from source row into target row)
$map = @{
    4 = 10
    5 = 11
    6 = 17
    7 = 18
    8 = 14
    9 = 15
    10 = 4
    11 = 5
    12 = 6
    13 = 22
    14 = 23
    15 = 24
    16 = 25
    17 = 16
    18 = 28
    19 = 29
    20 = 30
    21 = 26
    22 = 27
    23 = 31
    24 = 32
    25 = 33
    26 = 21
    27 = 12
    28 = 13
    29 = 19
    30 = 20
    31 = 7
    32 = 8
    33 = 9
}

foreach ($r in 4..33) {
    $src = $map[$r]
    $srcVals = $orig[$src]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $srcVals[$c]
    }
}

$wb.Save()
